$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "29.628.92"
Set-TextValue $ws.Range("E2") "  -0.56%  "
Set-TextValue $ws.Range("D3") "1.922.25"
Set-TextValue $ws.Range("E3") "  -1.41%  "
Set-TextValue $ws.Range("D4") "0.9966"
Set-TextValue $ws.Range("E4") "  -0.41%  "
Set-TextValue $ws.Range("D5") "334.48"
Set-TextValue $ws.Range("E5") "  -2.04%  "
Set-TextValue $ws.Range("D6") "0.9974"
Set-TextValue $ws.Range("E6") "  -0.31%  "
Set-TextValue $ws.Range("D7") "0.4659"
Set-TextValue $ws.Range("E7") "  -2.78%  "
Set-TextValue $ws.Range("D8") "0.4163"
Set-TextValue $ws.Range("E8") "  +0.70%  "
Set-TextValue $ws.Range("D9") "48.34"
Set-TextValue $ws.Range("E9") "  +1.05%  "
Set-TextValue $ws.Range("D10") "0.08056"
Set-TextValue $ws.Range("E10") "  -2.15%  "
Set-TextValue $ws.Range("E11") "  -1.41%  "
Set-TextValue $ws.Range("D12") "22.37"
Set-TextValue $ws.Range("E12") "  -1.73%  "
Set-TextValue $ws.Range("D13") "1.918.97"
Set-TextValue $ws.Range("E13") "  -2.38%  "
Set-TextValue $ws.Range("D14") "6.007"
Set-TextValue $ws.Range("E14") "  -2.39%  "
Set-TextValue $ws.Range("D15") "7.180"
Set-TextValue $ws.Range("E15") "  -2.83%  "
Set-TextValue $ws.Range("D16") "89.76"
Set-TextValue $ws.Range("E16") "  -2.32%  "
Set-TextValue $ws.Range("D17") "0.9967"
Set-TextValue $ws.Range("E17") "  -0.49%  "
Set-TextValue $ws.Range("E18") "  -2.14%  "
Set-TextValue $ws.Range("E19") "  -1.21%  "
Set-TextValue $ws.Range("D20") "17.83"
Set-TextValue $ws.Range("E20") "  -0.95%  "
Set-TextValue $ws.Range("D21") "0.9966"
Set-TextValue $ws.Range("E21") "  -0.37%  "
Set-TextValue $ws.Range("D22") "29.565.25"
Set-TextValue $ws.Range("E22") "  -0.64%  "
Set-TextValue $ws.Range("D23") "5.545"
Set-TextValue $ws.Range("E23") "  -0.73%  "
Set-TextValue $ws.Range("D24") "11.48"
Set-TextValue $ws.Range("E24") "  +1.98%  "
Set-TextValue $ws.Range("D25") "2.202"
Set-TextValue $ws.Range("E25") "  -3.91%  "
Set-TextValue $ws.Range("D26") "2.131.64"
Set-TextValue $ws.Range("E26") "  -2.55%  "
Set-TextValue $ws.Range("D27") "156.73"
Set-TextValue $ws.Range("E27") "  -2.77%  "
Set-TextValue $ws.Range("D28") "19.94"
Set-TextValue $ws.Range("E28") "  -1.43%  "
Set-TextValue $ws.Range("D29") "2.165"
Set-TextValue $ws.Range("E29") "  -0.36%  "
Set-TextValue $ws.Range("D30") "5.672"
Set-TextValue $ws.Range("E30") "  +0.22%  "
Set-TextValue $ws.Range("D31") "117.54"
Set-TextValue $ws.Range("E31") "  -4.37%  "
Set-TextValue $ws.Range("E32") "  +3.67%  "
Set-TextValue $ws.Range("D33") "0.09465"
Set-TextValue $ws.Range("E33") "  -2.03%  "
Set-TextValue $ws.Range("E34") "  -2.26%  "
Set-TextValue $ws.Range("D35") "5.458"
Set-TextValue $ws.Range("E35") "  -0.78%  "
Set-TextValue $ws.Range("D36") "3.532"
Set-TextValue $ws.Range("E36") "  -4.09%  "
Set-TextValue $ws.Range("D37") "0.06141"
Set-TextValue $ws.Range("E37") "  -1.89%  "
Set-TextValue $ws.Range("E38") "  -2.08%  "
Set-TextValue $ws.Range("D39") "8.462"
Set-TextValue $ws.Range("E39") "  -0.54%  "
Set-TextValue $ws.Range("D40") "1.178"
Set-TextValue $ws.Range("E40") "  -0.74%  "
Set-TextValue $ws.Range("D41") "0.5923"
Set-TextValue $ws.Range("E41") "  -2.61%  "
Set-TextValue $ws.Range("D42") "0.9968"
Set-TextValue $ws.Range("E42") "  -0.35%  "
Set-TextValue $ws.Range("E43") "  -4.20%  "
Set-TextValue $ws.Range("D44") "0.1839"
Set-TextValue $ws.Range("E44") "  -2.94%  "
Set-TextValue $ws.Range("D45") "2.372"
Set-TextValue $ws.Range("E45") "  -0.98%  "
Set-TextValue $ws.Range("E46") "  -2.55%  "
Set-TextValue $ws.Range("D47") "0.07543"
Set-TextValue $ws.Range("E47") "  +1.58%  "
Set-TextValue $ws.Range("D48") "0.5597"
Set-TextValue $ws.Range("E48") "  -1.89%  "
Set-TextValue $ws.Range("D49") "12.18"
Set-TextValue $ws.Range("D50") "1.941"
Set-TextValue $ws.Range("E50") "  -2.54%  "
Set-TextValue $ws.Range("D51") "112.86"
Set-TextValue $ws.Range("E51") "  -0.18%  "
